# hoursLog.xlsx update:
#  - created another object (row 7 C/D/E), can modify the object (row 8),
#    add outlets (row 9), plus placeholder shared-formula rows 10-20
#  - widened column E so the longer notes are readable
#  - left the selection on E10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (closest achievable value to 99.6640625 chars given this
# engine's pixel-rounded column-width model)
$ws.Columns.Item(5).ColumnWidth = 98.8

# --- Row 7: finish the entry that already had a start time (B7) ---
$ws.Range("C7").Value2 = 0.74305555555555547
$ws.Range("C7").NumberFormat = "h:mm"

# --- Row 8: new time-in/time-out entry ---
$ws.Range("B8").Value2 = 0.75
$ws.Range("B8").NumberFormat = "h:mm"
$ws.Range("C8").Value2 = 0.79861111111111116
$ws.Range("C8").NumberFormat = "h:mm"
$ws.Rows.Item(8).RowHeight = 30

# --- Row 9: new time-in/time-out entry ---
$ws.Range("B9").Value2 = 0.34027777777777773
$ws.Range("B9").NumberFormat = "h:mm"
$ws.Range("C9").Value2 = 0.41666666666666669
$ws.Range("C9").NumberFormat = "h:mm"

# Duration formula, filled as one shared formula across D7:D20 (matches
# the workbook's existing pattern of a "buffer" of blank duration rows)
$ws.Range("D7:D20").Formula = "=C7-B7"
$ws.Range("D7:D20").NumberFormat = "h:mm"

# Notes column for the three new log entries
$ws.Range("E7").Value2 = "able to build a project and build objects within that project, using the IDE"
$ws.Range("E8").Value2 = "creating help patches for new object, how does modifying a patch work in the ecosystem. Update comes from the .xml file in the doc folder"
$ws.Range("E9").Value2 = "able to create a new inlet, now working on inlets, arguments, and attributes"

# Leave the active selection on E10, like the author's last click
$ws.Range("E10").Select() | Out-Null
